# Auto-update draw results: append the 2025-10-29 Pick 3 row to the
# "Results" sheet (row 43), mirroring the existing rows 2-42 layout:
#   A: draw date            B: game name   C: phase code
#   D: result (digits)      E: ISO-8601 insert timestamp
#
# All five columns are stored as literal text in the source file (even
# though some values look numeric/date-like), so values that could be
# auto-converted by Excel's type inference (dates, pure-digit strings)
# are entered with a leading apostrophe to force text entry. The
# Style reset afterwards clears the "number stored as text" quote-prefix
# formatting Excel applies automatically, so the new row's cells keep
# the same (default/unstyled) look as every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43

$ws.Range("A$row").Value = "'2025-10-29"
$ws.Range("B$row").Value = "Pick 3"
$ws.Range("C$row").Value = "'251029"
$ws.Range("D$row").Value = "5-8-4"
$ws.Range("E$row").Value = "'2025-10-29T21:40:32.995+04:00"

# Strip the auto-applied "text-as-number" styling so the new cells match
# the unstyled look of the rest of the table.
$ws.Range("A$row`:E$row").Style = "Normal"
